$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.1375548059379061
$ws.Range("E2").Value = 9.879873151939263
$ws.Range("F2").Value = 31.48854667218161
